$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Change cell A1 value from "c-membership-type" to "Elements"
$ws.Range("A1").Value = "Elements"

# Change selection to A2
$ws.Range("A2").Select()
